$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B for "distance_miles"
$ws.Range("B1").EntireColumn.Insert()

# Set header for new column
$ws.Range("B1").Value = "distance_miles"

# New route data (route_id, distance_miles, origin, destination, num_flights,
# avg_revenue, avg_cost, avg_margin, avg_profit_pct, delay_count, avg_delay_minutes)
$data = @(
    @("RT4447", 2475, "JFK", "LAX", 13, 42268.72, 48264.62, -5995.9, -15.79, 3, 36.67),
    @("RT5358", 1846, "SFO", "ORD", 10, 37228.61, 41804.87, -4576.26, -12.69, 4, 67),
    @("RT6259", 3625, "JFK", "CDG", 10, 114673.23, 99245.39, 15427.84, 12.26, 3, 49),
    @("RT6474", 1258, "BOS", "MIA", 9, 45638.98, 32079.6, 13559.38, 29.13, 1, 62),
    @("RT6565", 1024, "SEA", "DEN", 9, 40214.39, 28622.01, 11592.38, 28.67, 0, 0),
    @("RT9374", 8446, "SFO", "SIN", 9, 213257.46, 171190.24, 42067.22, 18.47, 1, 107),
    @("RT7502", 732, "ATL", "DFW", 7, 34668.01, 26577.93, 8090.07, 18.37, 1, 120)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
